$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.472.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.52%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.325.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.29%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'310.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'108.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +4.19%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E8").Value = "'  +0.03%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.614"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +2.20%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'40.84"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.0916"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.22%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +1.58%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -0.79%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -0.62%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'15.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'2.679.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.24%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.329.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.75%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'43.360.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.37%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +1.20%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0000107"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.16%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'13.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -5.71%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.18%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'3.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.81%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'268.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.57%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.69%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.20%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'7.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +6.69%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'11.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.34%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -2.58%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'38.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.97%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +1.00%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'167.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.19%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.0886"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.75%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'2.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +9.47%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.78%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'4.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.93%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -2.00%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +3.95%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'NEARProtocol"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'3.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.15%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'LidoDAOToken"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'2.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +5.79%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +8.04%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'104.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +13.98%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'71.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +3.49%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +2.82%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'13.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +9.05%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +0.07%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'113.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.36%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.659.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.06%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'5.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.70%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'8.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'1.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +9.81%  "
$ws.Range("E51").Style = "Normal"
